# Generate Report for Handback
# This script updates the localization-status workbook to reflect that the
# handback has completed for both target locales (zh-cn and de-de):
#  - Overview sheet status cells flip from "Ready for handoff" to
#    "Handed back: in sync with en-US"
#  - Each locale sheet gets its "Latest Target File" / "Latest Handback
#    File" / "Latest Handback DateTime" columns populated, with the target
#    file cell turned into a hyperlink back to the source doc (just like
#    column A already is)
#  - A few columns are widened so the new long values are readable

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Literal values re-used throughout
# ---------------------------------------------------------------------
$statusHandedBack = "Handed back: in sync with en-US"

$name1812 = "1812dcea-ea10-4e77-aa74-f5f579ef231e.md"
$nameA216 = "a2160a50-b543-48d5-b194-3f1d31dfe14b.md"

$url1812 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e5fafeb7896e9a20ae446ba1c6c9cdec4aa0175c/e2e/1812dcea-ea10-4e77-aa74-f5f579ef231e.md"
$urlA216 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e5fafeb7896e9a20ae446ba1c6c9cdec4aa0175c/e2e/a2160a50-b543-48d5-b194-3f1d31dfe14b.md"

$zhcnTarget1812 = "1812dcea-ea10-4e77-aa74-f5f579ef231e.65d2812003b98068b846c75440c6a993de87736a.zh-cn.xlf"
$zhcnTargetA216 = "a2160a50-b543-48d5-b194-3f1d31dfe14b.0dffd52158bacd2672cc337495db96016c9628a0.zh-cn.xlf"
$dedeTarget1812 = "1812dcea-ea10-4e77-aa74-f5f579ef231e.65d2812003b98068b846c75440c6a993de87736a.de-de.xlf"
$dedeTargetA216 = "a2160a50-b543-48d5-b194-3f1d31dfe14b.0dffd52158bacd2672cc337495db96016c9628a0.de-de.xlf"

$zhcnHandbackDate = "2016-09-07 15:11:30"
$dedeHandbackDate = "2016-09-07 15:11:55"

$wideWidth = 29.15   # renders as the ~29.98-char column width used for Status columns
$fullWidth = 39.15   # renders as an exact 40-char column width

# ---------------------------------------------------------------------
# Overview sheet: flip the "Ready for handoff" status to "Handed back"
# for both locale rows/columns
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $statusHandedBack
$wsOverview.Range("F2").Value = $statusHandedBack
$wsOverview.Range("E3").Value = $statusHandedBack
$wsOverview.Range("F3").Value = $statusHandedBack

$wsOverview.Columns.Item(5).ColumnWidth = $wideWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideWidth

# ---------------------------------------------------------------------
# Helper-style block applied to both locale sheets (zh-cn / de-de): fill
# in Latest Target File (I), Latest Handback File (J) and Latest Handback
# DateTime (K) for both data rows, and rebuild the hyperlinks collection
# so the existing Source File Name links (A2/A3) keep working and the new
# Latest Target File links (I2/I3) point at the same source document.
# ---------------------------------------------------------------------

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("J2").Value = $zhcnTarget1812
$wsZh.Range("K2").Value = $zhcnHandbackDate
$wsZh.Range("J3").Value = $zhcnTargetA216
$wsZh.Range("K3").Value = $zhcnHandbackDate

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $url1812, $null, $null, $name1812) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $url1812, $null, $null, $name1812) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $urlA216, $null, $null, $nameA216) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $urlA216, $null, $null, $nameA216) | Out-Null

$wsZh.Columns.Item(3).ColumnWidth = $wideWidth
$wsZh.Columns.Item(9).ColumnWidth = $fullWidth
$wsZh.Columns.Item(10).ColumnWidth = $fullWidth

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("J2").Value = $dedeTarget1812
$wsDe.Range("K2").Value = $dedeHandbackDate
$wsDe.Range("J3").Value = $dedeTargetA216
$wsDe.Range("K3").Value = $dedeHandbackDate

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $url1812, $null, $null, $name1812) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $url1812, $null, $null, $name1812) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $urlA216, $null, $null, $nameA216) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $urlA216, $null, $null, $nameA216) | Out-Null

$wsDe.Columns.Item(3).ColumnWidth = $wideWidth
$wsDe.Columns.Item(9).ColumnWidth = $fullWidth
$wsDe.Columns.Item(10).ColumnWidth = $fullWidth
